$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data inputs -------------------------------------------------
# Contratos Negociados (number of contracts traded)
$ws.Range("B3").Value = 44
# Lucro/Prejuizo Bruto (gross profit/loss)
$ws.Range("B5").Value = 12

# New daily row (18/06/2018) in the "Ultimos 21 Pregoes" table
$ws.Range("J7").Value = "18/06/2018"
$ws.Range("K7").Value = 14
$ws.Range("L7").Value = 22.58

# --- Formula changes ---------------------------------------------------
# Lucro/Prejuizo Liquido now also nets out the new IRRF (1%) withholding
$ws.Range("E4").Formula = "=D27-D28"

# New IRRF (1%) line under "Ajuste DayTrade"
$ws.Range("C28").Value = "IRRF (1%)"
$ws.Range("C28").HorizontalAlignment = -4152
$ws.Range("D28").Formula = "=ROUND(IF(D27>0,D27*0.01,0),2)"
$ws.Range("D28").HorizontalAlignment = -4108

# --- Selection / cursor -------------------------------------------------
$ws.Range("C5").Select()
